$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the duplicate bold title paragraph near the end of the document
#    ("Play Big Ghoulies for Free - Review of Big Ghoulies Slot Game"). It is
#    the second-to-last paragraph in the (still unmodified) document.
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$dupePara = $d.Paragraphs($count - 1)
$dupePara.Range.Delete() | Out-Null

# ---------------------------------------------------------------------------
# 2) Replace the italic closing paragraph's text (now the last paragraph)
#    with the new image-prompt copy. We locate it with Find, then assign
#    straight-quote-safe text directly to the Range (Find.Execute's
#    replacement argument would smart-quote the literal double quotes).
# ---------------------------------------------------------------------------
$count2 = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($count2)
$imgRange = $lastPara.Range
$imgRange.Find.Execute("Play Big Ghoulies for free and read our review of this horror-themed slot machine. Simple gameplay, impressive graphics, dark atmosphere, and more.") | Out-Null
$imgRange.Text = 'Create a feature image for the game "Big Ghoulies" that fits the following criteria: Design a cartoon-style image featuring a happy Maya warrior with glasses. The warrior should be shown in a spooky setting, surrounded by symbols of horror such as ghosts, pumpkins, and skeletons. The warrior should be holding a treasure chest filled with gold coins and gems, representing the theme of the game. The image should use a dark color palette, with contrasting pops of bright colors to make the treasure stand out. The overall composition should be balanced and visually appealing, with attention paid to detail and texture. The Maya warrior should be depicted with a friendly and approachable expression, inviting players to join in on the spooky fun of Big Ghoulies.'

# ---------------------------------------------------------------------------
# 3) Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph at the top of the document.
# ---------------------------------------------------------------------------
$title = $d.Paragraphs(1)
$title.Range.InsertParagraphAfter() | Out-Null
$metaPara = $d.Paragraphs(2)

$metaXml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
           "<w:r/>" +
           "<w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>" +
           "<w:r><w:t>: Play Big Ghoulies for free and read our review of this horror-themed slot machine. Simple gameplay, impressive graphics, dark atmosphere, and more.</w:t></w:r>" +
           "</w:p>"
$metaPara.Range.InsertXML($metaXml) | Out-Null

Write-Output "done"
